$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert the newly-surfaced states as new rows, in the order they appear in
# the final table. Inserting repeatedly at the same anchor row pushes
# earlier insertions further down, so each group is written back-to-front
# (last name of the group first) so it ends up forward-ordered in the sheet.
function Insert-Rows-At {
    param($row, $names)
    for ($i = $names.Count - 1; $i -ge 0; $i--) {
        $ws.Rows.Item($row).Insert()
        $ws.Cells.Item($row, 1).Value = $names[$i]
    }
}

Insert-Rows-At 2  @("Minas Gerais")
Insert-Rows-At 14 @("Goiás", "Mato Grosso")
Insert-Rows-At 17 @("Amazonas")
Insert-Rows-At 20 @("Mato Grosso do Sul")
Insert-Rows-At 23 @("Distrito Federal", "Rondônia", "Acre")
Insert-Rows-At 27 @("Tocantins", "Roraima")

# Final (count, proportion) pairs for every data row 2..29, in order. The
# interviewee counts for pre-existing states are unchanged; only the
# proportion (now computed against the larger 43948 total) moves. New rows
# get both their count and proportion filled in here too.
$values = @(
    @(3505, 0.0797533448621098),
    @(2923, 0.0665104214071175),
    @(2402, 0.0546555019568581),
    @(2353, 0.0535405479202694),
    @(2278, 0.0518339856193683),
    @(2063, 0.0469418403567853),
    @(2035, 0.0463047237644489),
    @(1997, 0.0454400655319924),
    @(1808, 0.0411395285337217),
    @(1795, 0.0408437244015655),
    @(1757, 0.0399790661691089),
    @(1754, 0.0399108036770729),
    @(1691, 0.038477291344316),
    @(1537, 0.0349731500864658),
    @(1535, 0.0349276417584418),
    @(1500, 0.0341312460180213),
    @(1434, 0.0326294711932284),
    @(1425, 0.0324246837171202),
    @(1287, 0.0292846090834623),
    @(1218, 0.0277145717666333),
    @(1215, 0.0276463092745973),
    @(980, 0.0222990807317739),
    @(737, 0.0167698188768545),
    @(731, 0.0166332938927824),
    @(709, 0.0161327022845181),
    @(662, 0.0150632565759534),
    @(617, 0.0140393191954128),
    @(43948, 1)
)

$r = 2
foreach ($pair in $values) {
    $ws.Cells.Item($r, 2).Value = $pair[0]
    $ws.Cells.Item($r, 3).Value = $pair[1]
    $r = $r + 1
}
